$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Edit A2's text in place: "hand" -> "left/right hand"
$ws.Range("A2").Value = "I have/has the doc/pen in my left/right hand"

# Clear the contents of A3 and A4 (keep their existing style/formatting)
$ws.Range("A3").ClearContents()
$ws.Range("A4").ClearContents()

# Move the active selection from A6 to A4
$ws.Range("A4").Select() | Out-Null
